$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set values 1 in G2:G7, H2:H7, and I2:I6 (I7 already has content)
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 7).Value = 1   # column G
    $ws.Cells.Item($r, 8).Value = 1   # column H
}
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 9).Value = 1   # column I
}

# Update the selection to match the new active cell/range
$ws.Activate()
$ws.Range("H8").Select()
